$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "El archivo comprimido..." paragraph: the sentence was split across two
#    runs ("...para e" / "l funcionamiento...") -- reunite it into a single
#    run by replacing the full sentence text with itself.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "El archivo comprimido contiene todo lo necesario para el funcionamiento del programa.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El archivo comprimido contiene todo lo necesario para el funcionamiento del programa.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop the "El resultado final...: / resultado.flg" paragraph completely
#    (bold/italic paragraph + its following blank paragraph).
# ---------------------------------------------------------------------------
$resultParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "El resultado final del proceso*") {
        $resultParaIndex = $i
        break
    }
}
if ($resultParaIndex -eq 0) {
    Write-Output "WARNING: 'El resultado final...' paragraph not found"
} else {
    # Delete the paragraph itself (removes text + mark), then delete the
    # blank paragraph that used to sit right after it.
    $d.Paragraphs.Item($resultParaIndex).Range.Delete() | Out-Null
    $d.Paragraphs.Item($resultParaIndex).Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 3) "NERC_GIL.bat ARCHIVO" paragraph gains a third run: "_ENTRADA ARCHIVO_SALIDA"
# ---------------------------------------------------------------------------
$nercParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "NERC_GIL.bat ARCHIVO*") {
        $nercParaIndex = $i
        break
    }
}
if ($nercParaIndex -eq 0) {
    Write-Output "WARNING: 'NERC_GIL.bat ARCHIVO' paragraph not found"
} else {
    $nercPara = $d.Paragraphs.Item($nercParaIndex)
    $insertPos = $d.Range($nercPara.Range.End - 1, $nercPara.Range.End - 1)
    $insertPos.InsertAfter("_ENTRADA ARCHIVO_SALIDA") | Out-Null
}

# ---------------------------------------------------------------------------
# 4) "En donde ARCHIVO es el nombre..." paragraph: rewrite its text, drop the
#    trailing picture, and keep the "_GoBack" bookmark anchored at the very
#    end of the paragraph (just before the paragraph mark).
# ---------------------------------------------------------------------------
$dondeParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "En donde ARCHIVO es el nombre*") {
        $dondeParaIndex = $i
        break
    }
}
if ($dondeParaIndex -eq 0) {
    Write-Output "WARNING: 'En donde ARCHIVO es el nombre...' paragraph not found"
} else {
    $dondePara = $d.Paragraphs.Item($dondeParaIndex)
    # Removing the whole paragraph (text + mark) also discards the picture
    # and the bookmark that lived inside it.
    $dondePara.Range.Delete() | Out-Null

    $afterPara = $d.Paragraphs.Item($dondeParaIndex)
    $insPt = $d.Range($afterPara.Range.Start, $afterPara.Range.Start)
    $xmlFrag = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00854193" w:rsidRPr="00854193" w:rsidRDefault="00F92086" w:rsidP="00F92086"><w:r><w:t>En donde ARCHIVO</w:t></w:r><w:r><w:t>_ENTRADA</w:t></w:r><w:r><w:t xml:space="preserve"> es el nombre del archivo que se busca analizar, </w:t></w:r><w:r><w:t>y ARCHIVO_SALIDA es el archivo que se desea crear con el resultado</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
"@
    $insPt.InsertXML($xmlFrag) | Out-Null

    # InsertXML above added one extra blank paragraph; remove it again so
    # the trailing blank paragraph count matches the original document.
    $d.Paragraphs.Item($dondeParaIndex + 1).Range.Delete() | Out-Null
}
